# M1_OriginMonitor // Simulation temp-pass edit
# Removes the redundant "Valve" label/prefix from the ADV/RET cylinder
# shapes on slide 3 ("Flow" / 직사각형 & 화살표 shapes), leaving the
# plain "ADV"/"RET" captions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape 2 ("직사각형 55"): two paragraphs "Valve" / "ADV" -> single "ADV"
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 6).Delete()

# Shape 3 ("직사각형 66"): two paragraphs "Valve" / "RET" -> single "RET"
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 6).Delete()

# Shape 4 ("화살표: 오각형 79"): "ADV[ValveADV ~ ValveADV]" -> "ADV[ADV ~ ADV]"
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Characters(6, 5).Delete()
$tr.Characters(12, 5).Delete()

# Shape 5 ("화살표: 오각형 80"): "RET[ValveRET ~ ValveRET]" -> "RET[RET ~ RET]"
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange
$tr.Characters(6, 5).Delete()
$tr.Characters(12, 5).Delete()
